$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (old row 2, "H 72"); all subsequent rows shift up by one.
$ws.Rows.Item(2).Delete()

# Apply the individual cell corrections (missing-value imputations / clears) that differ
# from a pure shift-up of the original data.
$ws.Range("F4").Value = 0.70909
$ws.Range("F8").ClearContents()
$ws.Range("F9").Value = 0.71194
$ws.Range("C11").Value = 15.5
$ws.Range("F11").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("F20").Value = 0.7106
$ws.Range("F22").ClearContents()
$ws.Range("F23").Value = 0.70931
$ws.Range("F25").ClearContents()
$ws.Range("C26").Value = 12
$ws.Range("C27").ClearContents()
$ws.Range("F28").Value = 0.70963
$ws.Range("C31").Value = 10.7
$ws.Range("C32").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("C35").Value = 11.3
$ws.Range("C36").ClearContents()
$ws.Range("F36").Value = 0.71087
$ws.Range("C37").Value = 12.1
$ws.Range("C38").ClearContents()
$ws.Range("F38").ClearContents()
$ws.Range("F41").Value = 0.71115
$ws.Range("F43").ClearContents()
$ws.Range("C45").Value = 11.7
$ws.Range("C46").ClearContents()
$ws.Range("F52").Value = 0.70948
$ws.Range("C53").Value = 10.5
$ws.Range("F54").ClearContents()
$ws.Range("C56").ClearContents()
